$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cesar Jimenez (row 4): add his e-mail (with mailto hyperlink), cell number
# and a narrower Provincias list.
$ws.Range("B4").Value = "cesarjjxd@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:cesarjjxd@gmail.com")
$ws.Range("B4").Style = $ws.Range("B2").Style

$ws.Range("C4").Value = "8529-6827"

# Fernanda Murillo (row 3): new cell number.
$ws.Range("C3").Value = "8598-6048"

$ws.Range("D4").Value = "San José,Heredia,Cartago"
$ws.Range("D3").Value = "San José,Heredia,Cartago,Alajuela,Puntarenas"

# Maynor Martinez (row 2): expanded province list (same text as row 3's).
$ws.Range("D2").Value = "San José,Heredia,Cartago,Alajuela,Puntarenas"

# Column D is now wider to fit the longer province lists (~48.57 characters).
$ws.Columns.Item(4).ColumnWidth = 47.67

# Selection moved to D2 in the saved view.
$ws.Range("D2").Select() | Out-Null
